# "fixing data and model"
# Adds four new letter/punctuation entries (p-sofit "tsadi", "tsadi geresh",
# open paren, close paren) as new rows at the bottom of the letters column,
# and updates the view (zoom/scroll/selection) to reflect the new scroll
# position in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the four new values to column A, rows 35-38.
$ws.Cells.Item(35, 1).Value = "ץ"
$ws.Cells.Item(36, 1).Value = "ץ'"
$ws.Cells.Item(37, 1).Value = "("
$ws.Cells.Item(38, 1).Value = ")"

# Update the view: zoom in to 200%, select A39 (first empty cell below the
# new data) and scroll so row 28 is at the top of the window.
$win = $excel.ActiveWindow
$win.Zoom = 200

$ws.Range("A39").Select() | Out-Null

$win.ScrollRow = 28
$win.ScrollColumn = 1
